$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 ("Experimental" row): set Value to literal text "false" (not boolean TRUE/FALSE).
# A direct Range.Value assignment of "false" is auto-coerced to a Boolean cell by
# the engine's type inference (like typing into Excel), so instead we write it as
# a text-returning formula and then convert that formula to a static value via
# Copy + PasteSpecial(xlPasteValues) - this preserves the literal string type.
$ws.Range("B7").Formula = '="false"'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# B8 ("Date" row): update the ValueSet date/time string.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# B17 ("Description" row): add the ValueSet description text.
$ws.Range("B17").Value = "Methods for comparing current values to established baselines"
